# Updated for TFS 6246:
#  - AT_Role_Access: add rows for job code WPSM13 (Sr Analyst, Functional) and
#    append the corresponding SQL INSERT lines (column K), fixing the missing
#    trailing comma on the previous (WPPM50) line.
#  - Revision_History: log the change as revision 6 / TFS# 6246.

$wb = $excel.ActiveWorkbook

# --- AT_Role_Access sheet ---
$wsRA = $wb.Worksheets.Item("AT_Role_Access")

# Row 13: WPSM13 / Sr Analyst, Functional -> RoleId 101 (CoachingAdmin)
$wsRA.Range("A13").Value = "WPSM13"
$wsRA.Range("B13").Value = "Sr Analyst, Functional"
$wsRA.Range("C13").Value = 101
$wsRA.Range("D13").Value = "CoachingAdmin"
$wsRA.Range("E13").Value = 0
$wsRA.Range("F13").Value = 1

# Row 14: WPSM13 / Sr Analyst, Functional -> RoleId 103 (WarningAdmin)
$wsRA.Range("A14").Value = "WPSM13"
$wsRA.Range("B14").Value = "Sr Analyst, Functional"
$wsRA.Range("C14").Value = 103
$wsRA.Range("D14").Value = "WarningAdmin"
$wsRA.Range("E14").Value = 0
$wsRA.Range("F14").Value = 1

# Column K holds the raw SQL script text; fix the trailing comma on the
# previous last line and append the two new INSERT value rows for WPSM13.
$wsRA.Range("K19").Value = "           ('WPPM50','Manager, Program',102,'CoachingUser',1,1),"
$wsRA.Range("K20").Value = "           ('WPSM13','Sr Analyst, Functional',101,'CoachingAdmin',0,1),"
$wsRA.Range("K21").Value = "           ('WPSM13','Sr Analyst, Functional',103,'WarningAdmin',0,1)"

$null = $wsRA.Range("M14").Select()

# --- Revision_History sheet ---
$wsRH = $wb.Worksheets.Item("Revision_History")

$wsRH.Range("A8").Value = 6
$wsRH.Range("B8").Value = 42836
# Reuse the date format from the row above instead of setting a fresh
# NumberFormat string (which would create a redundant style entry).
$null = $wsRH.Range("B7").Copy()
$null = $wsRH.Range("B8").PasteSpecial(-4122)
$wsRH.Range("C8").Value = "Susmitha Palacherla"
$wsRH.Range("D8").Value = 6246
$wsRH.Range("E8").Value = "Added rows for job code WPSM13 (Mark Hackman)AT_Role_Access tab"

$null = $wsRH.Activate()
$null = $wsRH.Range("E8").Select()
